# Apply the "Converted timing to frames" edit to cond_pm1.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: swap the shared-string labels of A1 and B1 ---
# Before: A1 = "angle_diff", B1 = "SOA"
# After:  A1 = "SOA",        B1 = "angle_diff"
$a1 = $ws.Range("A1").Value()
$b1 = $ws.Range("B1").Value()
$ws.Range("A1").Value = $b1
$ws.Range("B1").Value = $a1

# --- Data rows: convert SOA (column B, ms) into frames (column A) ---
# and move the old "angle bucket" flag (column A) into column B.
# The data is laid out in 8 contiguous blocks of 24 rows each:
#   rows 2-25    : angleFlag=0, SOA=0    -> frames=0,  angleFlag=0
#   rows 26-49   : angleFlag=0, SOA=20   -> frames=2,  angleFlag=0
#   rows 50-73   : angleFlag=0, SOA=60   -> frames=6,  angleFlag=0
#   rows 74-97   : angleFlag=0, SOA=120  -> frames=12, angleFlag=0
#   rows 98-121  : angleFlag=1, SOA=0    -> frames=0,  angleFlag=1
#   rows 122-145 : angleFlag=1, SOA=20   -> frames=2,  angleFlag=1
#   rows 146-169 : angleFlag=1, SOA=60   -> frames=6,  angleFlag=1
#   rows 170-193 : angleFlag=1, SOA=120  -> frames=12, angleFlag=1

$soaToFrames = @{0 = 0; 20 = 2; 60 = 6; 120 = 12}
$soaOrder = @(0, 20, 60, 120)
$blockSize = 24
$firstRow = 2
$lastRow = 193

for ($row = $firstRow; $row -le $lastRow; $row += $blockSize) {
    $endRow = $row + $blockSize - 1
    $blockIndex = ($row - $firstRow) / $blockSize
    $angleFlag = [Math]::Floor($blockIndex / 4)
    $soaIndex = $blockIndex % 4
    $soa = $soaOrder[$soaIndex]
    $frames = $soaToFrames[$soa]

    $ws.Range("A$row`:A$endRow").Value = $frames
    $ws.Range("B$row`:B$endRow").Value = $angleFlag
}

# --- Selection: whole column A selected, no fixed active cell ---
$ws.Range("A1:A1048576").Select()
